# Generate Report for Handoff
# Adds a new handed-off file (e8b12011-844e-4b1f-976f-d39a2a1085d5.md) as a
# new row (row 3) on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commit = "132c1e790752cde19d009d30a96e45851424364a"
$newFile = "e8b12011-844e-4b1f-976f-d39a2a1085d5.md"
$newFileHash = "e8b12011-844e-4b1f-976f-d39a2a1085d5.11fea0f6b71c07ccbad2f3789b7aa96d69f933e8"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name, Path And Name, Extension,
# Publish URL, zh-cn, de-de, Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\" + $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 12:49:04"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile",
    "",
    "",
    "e2e\" + $newFile
)
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 15570276

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": columns Source File Name, File Extension,
# Status, Source Path, Priority, Content Duplicate, Latest Handoff File,
# Latest Handoff Datetime, Latest Target File, Latest Handback File,
# Latest Handback DateTime, Reference Tokens, To be localized,
# Dependency From, Has metadata, Error Detail
# ---------------------------------------------------------------------
$langs = @(
    @{ Sheet = "zh-cn"; Xlf = $newFileHash + ".zh-cn.xlf"; HandoffDate = "2016-08-30 12:48:54" },
    @{ Sheet = "de-de"; Xlf = $newFileHash + ".de-de.xlf"; HandoffDate = "2016-08-30 12:49:04" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $ws.Range("A3").Value = $newFile
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = "e2e"
    $ws.Range("E3").Value = "ht"
    $ws.Range("F3").Value = "False"
    $ws.Range("G3").Value = $lang.Xlf
    $ws.Range("H3").Value = $lang.HandoffDate
    $ws.Range("H3").NumberFormat = $dateFmt
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = "0001-01-01 00:00:00"
    $ws.Range("K3").NumberFormat = $dateFmt
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = "True"
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = "False"
    $ws.Range("P3").Value = ""

    $ws.Hyperlinks.Add(
        $ws.Range("A3"),
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile",
        "",
        "",
        $newFile
    )
    $ws.Range("A3").Font.Underline = $true
    $ws.Range("A3").Font.Color = 15570276

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P3"))
}
